$d = $word.ActiveDocument
$search = " la primera creación de nuestro primer código de Android utilizaremos un template pre cargado, el cual nos permitirá que de forma automática la aplicación cree una serie de códigos, los cuales podemos modificar mas adelante, pero eso nos es lo que nos interesa por el momento, "
$rng = $d.Content
$found = $rng.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("found: " + $found)
Write-Output ("rng start/end: " + $rng.Start + "/" + $rng.End)
Write-Output ("rng text: [" + $rng.Text + "]")
